$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, pushing old rows 13-21 down to 14-22.
# This naturally realigns all of the "label" cells in column A and
# reproduces the target row-height pattern (Excel copies formatting
# from the row above when inserting).
$ws.Rows.Item(13).Insert()

# The inserted row 13 picks up a stray empty A13 cell (formatted like
# the row above it) - the target layout has no A13 cell at all.
$ws.Range("A13").Clear()

# --- Row 10: Objetivos: value text (was incorrectly the teacher name) ---
$ws.Range("B10").Value = "Apresentar e analisar os conceitos básicos de monitoramento, suas aplicações práticas e as interfaces com os demais instrumentos de Política Ambiental."
$ws.Range("C10").Value = "Apresentar e analisar os conceitos básicos de monitoramento, suas aplicações práticas e as interfaces com os demais instrumentos de Política Ambiental."

# --- Row 13: Docentes responsáveis value (teacher name), no column A label ---
$ws.Range("B13").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Range("C13").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Range("B13").WrapText = $true
$ws.Range("B13").Font.Bold = $false

# --- Row 14: Programa resumido: value text ---
$ws.Range("B14").Value = "Monitoramento da qualidade ambiental."
$ws.Range("C14").Value = "Monitoramento da qualidade ambiental."

# --- Row 16: Programa: value text ---
$ws.Range("B16").Value = "Conceitos de qualidade ambiental, poluição, padrões de qualidade e de emissão. Conceito de monitoramento. Amostragem. Sistemas de monitoramento. Índices de qualidade. Monitoramento como parte integrante de sistema de gestão ambiental."
$ws.Range("C16").Value = "Conceitos de qualidade ambiental, poluição, padrões de qualidade e de emissão. Conceito de monitoramento. Amostragem. Sistemas de monitoramento. Índices de qualidade. Monitoramento como parte integrante de sistema de gestão ambiental."

# --- Row 19: Método: value text ---
$ws.Range("B19").Value = "Aula expositiva e exercícios dirigidos."
$ws.Range("C19").Value = "Aula expositiva e exercícios dirigidos."

# --- Row 20: Critério: value text ---
$ws.Range("B20").Value = "Média ponderada de exercícios e provas."
$ws.Range("C20").Value = "Média ponderada de exercícios e provas."

# --- Row 21: Norma de recuperação: value text ---
$ws.Range("B21").Value = "Prova única com nota igual ou superior a 5,0."
$ws.Range("C21").Value = "Prova única com nota igual ou superior a 5,0."

# --- Row 22: Bibliografia: value text (new row) ---
$biblio = "Porto, R.L.:. org.. Técnicas quantitativas para o gerenciamento de recursos hídricos. ABRH e Editora da Universidade. 1997.`nJames, A. ed., Mathematical models in water pollution control. John Wiley & Sons. 1989. `nMota, S.. Preservação e Conservação de Recursos Hídricos. ABES. 2a. edição. 1995.`nSewell, G.H. Administração e controle de qualidade ambiental. EPU. 1998.`nMacknight, A. Handbook of techniques for aquatic sediments sampling. McGraw Hill 1999. `nLoeb, A. Biological monitoring of aquatic systems. McGraw-Hill. 1998."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio

# --- Row heights: rows 15 and 21 become 60 (were 120), rows 17 and 22 become 120 ---
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120

# --- Column widths: separate column A's width definition from column B's ---
$ws.Columns.Item(2).ColumnWidth = 60.71
